$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "Different Service for different model"
$ws.Range("B11").Value = "Change name of  namespaces "
$ws.Range("B14").Value = "Extensions"
$ws.Range("B14").Font.Bold = $true
$ws.Range("B15").Value = "Model Object  - DTO Object"
$ws.Range("B16").Value = "Request Model Object  - Model Object"
$ws.Range("B17").Value = "Model Class extra functionality"
$ws.Range("B19").Value = "Implement and understand pagination classes"
$ws.Range("B21").Value = "Tweet"
$ws.Range("B22").Value = "Like count "
$ws.Range("B23").Value = "Dislike count"

$null = $ws.Range("F15").Select()
